$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits -----------------------------------------------------------

# G2: phone number re-typed without spaces, forced as text (leading apostrophe)
$ws.Range("G2").Value = "'55111234567"

# L2: CEP re-typed without the dash, forced as text; this cell also carries
# the underline font (copied from the header style) together with the
# quote-prefix forced-text flag.
$ws.Range("L2").Font.Underline = 2
$ws.Range("L2").Value = "'06010060"

# A2: Usuario value got bumped from "Karl_Otaner2" to "Karl_Otaner12"
$ws.Range("A2").Value = "Karl_Otaner12"

# --- Stray formatted-but-empty cells ---------------------------------------
# L3 and A7 ended up holding the header's underline style with no value.
$ws.Range("L3").Font.Underline = 2
$ws.Range("A7").Font.Underline = 2

# --- Column width -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.8333333333333

# --- Selection ---------------------------------------------------------------
[void]$ws.Range("A7").Select()
